# Trade #62 closed at 2026-02-18 00:23:32 - unknown UNKNOWN +0.000%
#
# This script applies the following changes to the live trading results
# workbook:
#   1. Summary sheet totals roll forward by one trade (capital, P&L, counts).
#   2. Strategy Status row for HighProbConvergence reflects the closed trade.
#   3. All Trades: trade #90 (HighProbConvergence) transitions OPEN -> CLOSED,
#      and a brand-new trade #119 (MarketMaking) is appended as OPEN.
#   4. HighProbConvergence sheet: trade #90 row updated to CLOSED.
#   5. MarketMaking sheet: new trade #119 row appended as OPEN.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.14
$wsSummary.Range("B4").Value = 0.25
$wsSummary.Range("B5").Value = 0.06
$wsSummary.Range("B6").Value = 90
$wsSummary.Range("B7").Value = 45
$wsSummary.Range("B9").Value = 50

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - HighProbConvergence row (row 3)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C3").Value = 100.26
$wsStatus.Range("D3").Value = 8
$wsStatus.Range("E3").Value = 0.27
$wsStatus.Range("F3").Value = 0.26
$wsStatus.Range("G3").Value = 87.5

# ---------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Trade #90 (HighProbConvergence) closes out on row 91
$wsAll.Range("G91").Value = 0.026531
$wsAll.Range("H91").Value = "CLOSED"
$wsAll.Range("I91").Value = 32.6541
$wsAll.Range("J91").Value = 0.01
$wsAll.Range("K91").Value = 100.26
$wsAll.Range("L91").Value = "early_exit"
$wsAll.Range("M91").Value = 0.09

# New trade #119 (MarketMaking) appended as row 120. Duplicate the last
# existing row (119) first so text-typed cells (dates/times/strings) keep
# their text type instead of being re-parsed (and possibly auto-converted
# to a date serial) as fresh literals, then overwrite only what differs.
$wsAll.Range("A119:Q119").Copy($wsAll.Range("A120:Q120"))
$wsAll.Range("A120").Value = 119
$wsAll.Range("C120").Value = "00:23:26"
$wsAll.Range("E120").Value = "UP"
$wsAll.Range("F120").Value = 0.02
# G120/L120 (no exit price / exit reason yet) must stay empty *text* cells
# (like the source row's blanks) rather than turn into numeric blanks -
# the Copy() above demotes them, so restore with a quote-prefix empty
# string and strip the resulting style back to Normal.
$wsAll.Range("G120").Value = "'"
$wsAll.Range("G120").Style = "Normal"
$wsAll.Range("L120").Value = "'"
$wsAll.Range("L120").Style = "Normal"

# ---------------------------------------------------------------------
# 4. HighProbConvergence sheet - trade #90 closes out on row 9
# ---------------------------------------------------------------------
$wsHPC = $wb.Worksheets.Item("HighProbConvergence")
$wsHPC.Range("G9").Value = 0.026531
$wsHPC.Range("H9").Value = "CLOSED"
$wsHPC.Range("I9").Value = 32.6541
$wsHPC.Range("J9").Value = 0.01
$wsHPC.Range("K9").Value = 100.26
$wsHPC.Range("P9").Value = "early_exit"
$wsHPC.Range("Q9").Value = 0.09

# ---------------------------------------------------------------------
# 5. MarketMaking sheet - new trade #119 appended as row 40
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Duplicate the last existing row (39) so text-typed cells keep their
# text type, then overwrite only what differs.
$wsMM.Range("A39:Q39").Copy($wsMM.Range("A40:Q40"))
$wsMM.Range("A40").Value = 119
$wsMM.Range("C40").Value = "00:23:26"
$wsMM.Range("E40").Value = "UP"
$wsMM.Range("F40").Value = 0.02
# G40 (no exit price yet) / P40 (no exit reason yet) must stay empty
# *text* cells rather than turn into numeric blanks - restore as above.
$wsMM.Range("G40").Value = "'"
$wsMM.Range("G40").Style = "Normal"
$wsMM.Range("P40").Value = "'"
$wsMM.Range("P40").Style = "Normal"
